$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.781.44'
$ws.Range("E2").Value = '  -1.05%  '

$ws.Range("D3").Value = '2.273.79'
$ws.Range("E3").Value = '  -1.13%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '249.89'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.23%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.634'
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.58%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '79.42'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +7.82%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.647'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '41.64'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.53%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0967'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '7.37'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.21%  '

$ws.Range("E13").Value = '  -1.22%  '

$ws.Range("D14").Value = '2.614.41'
$ws.Range("E14").Value = '  -0.88%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.12'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.871'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.40%  '

$ws.Range("D17").Value = '2.274.69'
$ws.Range("E17").Value = '  -0.89%  '

$ws.Range("D18").Value = '42.690.25'
$ws.Range("E18").Value = '  -1.19%  '

$ws.Range("D19").Value = '0.0₃0994'
$ws.Range("E19").Value = '  -1.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.21'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.66%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.00'
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.61%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '232.02'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.16'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.02%  '

$ws.Range("E24").Value = '  -3.09%  '

$ws.Range("E25").Value = '  -0.05%  '

$ws.Range("E26").Value = '  -5.08%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.33'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.44%  '

$ws.Range("E28").Value = '  +1.13%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '169.26'
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.47%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.78'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +6.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.81'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.59%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.0851'
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.122'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.17%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '30.49'
$ws.Range("D34").NumberFormat = "General"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -2.09%  '

$ws.Range("E35").Value = '  +0.87%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.56'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.19%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.76'
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.68%  '

$ws.Range("E38").Value = '  -3.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '13.46'
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.05%  '

$ws.Range("E40").Value = '  -4.20%  '

$ws.Range("E41").Value = '  -2.67%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '115.47'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.97%  '

$ws.Range("E43").Value = '  -2.70%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '61.46'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.87'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.21%  '

$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("B47").Value = 'BinanceUSD'
$ws.Range("C47").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.01%  '

$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.51'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -8.69%  '

$ws.Range("E49").Value = '  -4.67%  '

$ws.Range("E50").Value = '  -3.08%  '

$ws.Range("E51").Value = '  -2.33%  '
